$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the reporting-period dates on row 8 (B8:C8 period start/end, V8:W8 validation dates)
$ws.Cells.Item(8, 2).Value = Get-Date -Year 2021 -Month 7 -Day 1 -Hour 0 -Minute 0 -Second 0
$ws.Cells.Item(8, 3).Value = Get-Date -Year 2021 -Month 12 -Day 31 -Hour 0 -Minute 0 -Second 0
$ws.Cells.Item(8, 22).Value = Get-Date -Year 2022 -Month 1 -Day 10 -Hour 0 -Minute 0 -Second 0
$ws.Cells.Item(8, 23).Value = Get-Date -Year 2022 -Month 1 -Day 10 -Hour 0 -Minute 0 -Second 0

# Move the active view/selection to show the far right of the sheet
$ws.Application.ActiveWindow.ScrollColumn = 21
$ws.Range("W12").Select()
